# Update "Pais" sheet with refreshed COVID-19 country figures
# (data as of 19 Abril 2020, 11:22 instead of 10:52) and the resulting
# re-ranking of a handful of countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 11:22"

# Austria keeps rank 26 (row 22), figures refreshed
Set-Row 22 @("Austria", 14675, 4, 10501, 3731, 204, 0, 443)

# Ranks 41-46 (rows 37-42) reshuffle: Emiratos Arabes Unidos overtakes Chequia,
# and Singapur overtakes Indonesia/Filipinas
Set-Row 37 @("Emiratos Arabes Unidos", 6781, 479, 1286, 5454, 1, 4, 41)
Set-Row 38 @("Chequia", 6657, 51, 1235, 5241, 84, 0, 181)
Set-Row 39 @("Australia", 6606, 20, 4230, 2306, 48, 0, 70)
Set-Row 40 @("Singapur", 6588, 596, 740, 5837, 23, 0, 11)
Set-Row 41 @("Indonesia", 6575, 327, 686, 5307, 0, 47, 582)
Set-Row 42 @("Filipinas", 6259, 172, 572, 5278, 1, 12, 409)

# Rank 49 (row 45) Malasia refreshed
Set-Row 45 @("Malasia", 5389, 84, 3197, 2103, 46, 1, 89)

# Rank 54 (row 50) Finlandia refreshed
Set-Row 50 @("Finlandia", 3783, 102, 1700, 1993, 70, 0, 90)

# Rank 83 (row 79) Eslovaquia refreshed
Set-Row 79 @("Eslovaquia", 1161, 72, 229, 920, 10, 1, 12)

# Ranks 109-111 (rows 105-107) reshuffle: Estado de Palestina overtakes Malta/Taiwan
Set-Row 105 @("Estado de Palestina", 431, 13, 71, 358, 0, 0, 2)
Set-Row 106 @("Malta", 426, 0, 99, 324, 4, 0, 3)
Set-Row 107 @("Taiwan", 420, 22, 189, 225, 0, 0, 6)

# Ranks 171-172 (rows 167-168) reshuffle: Maldivas overtakes Mozambique
Set-Row 167 @("Maldivas", 36, 1, 16, 20, 0, 0, 0)
Set-Row 168 @("Mozambique", 35, 0, 4, 31, 0, 0, 0)
